# Increased angle after wave 10, up to 70
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column O ("Spread") used to be a flat 10 for waves 1-10 (rows 3-12).
# Ramp it up by 5 per wave instead: 15, 20, 25, ... 60.
# (Rows 13-24 already compute O(n-1)+1 via formula and will recalc
# automatically once O12 changes.)
$ws.Range("O3").Value = 15
$ws.Range("O4").Value = 20
$ws.Range("O5").Value = 25
$ws.Range("O6").Value = 30
$ws.Range("O7").Value = 35
$ws.Range("O8").Value = 40
$ws.Range("O9").Value = 45
$ws.Range("O10").Value = 50
$ws.Range("O11").Value = 55
$ws.Range("O12").Value = 60

# Leave the cursor where the author last left it while editing.
$ws.Range("D7").Select()
